# Add a per-column "max" row (13) above the existing overall "max" row (14),
# and rename the overall max label to "max of all".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: per-column max of abs, inserted just above the existing summary row ---
$ws.Range("R13").Value2 = "max"
$ws.Range("R13").Font.Bold = $true
$ws.Range("R13").HorizontalAlignment = -4108   # xlCenter

# Highlight the new max cells with a darker green fill (RGB 0,176,80 -> 0x00B050)
$ws.Range("S13:U13").Interior.Color = 5287936

$ws.Range("S13").FormulaArray = "=MAX(ABS(S3:S12))"
$ws.Range("T13").FormulaArray = "=MAX(ABS(T3:T12))"
$ws.Range("U13").FormulaArray = "=MAX(ABS(U3:U12))"

# --- Row 14: relabel the existing overall-max cell ---
$ws.Range("S14").Value2 = "max of all"

# --- Page setup / view tweaks ---
$ws.PageSetup.Orientation = 1   # xlPortrait
$ws.Range("U21").Select() | Out-Null
